# Insert a new data row for "Feria Lagunitas de Puerto Montt" - Pepino ensalada
# at row 104, shifting the existing rows 104:160 down to 105:161.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 104..160 down by one to make room for the new record.
$ws.Rows.Item(104).Insert()

# Populate the newly inserted row 104 with its data.
$ws.Cells.Item(104, 1).Value = 4
$ws.Cells.Item(104, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(104, 3).Value = "Los Lagos"
$ws.Cells.Item(104, 4).Value = 44488
$ws.Cells.Item(104, 5).Value = 10
$ws.Cells.Item(104, 6).Value = 100112043
$ws.Cells.Item(104, 7).Value = "Pepino ensalada"
$ws.Cells.Item(104, 8).Value = "Sin especificar"
$ws.Cells.Item(104, 9).Value = "Primera"
$ws.Cells.Item(104, 10).Value = 400
$ws.Cells.Item(104, 11).Value = 13000
$ws.Cells.Item(104, 12).Value = 15000
$ws.Cells.Item(104, 13).Value = 14000
$ws.Cells.Item(104, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(104, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(104, 16).Value = 233
$ws.Cells.Item(104, 17).Value = 60
$ws.Cells.Item(104, 18).Value = "Hortaliza"
